# SIS Personal Page & ProductDashboard changes
#
# Adds a batch of new "Personal Page" form columns (J:S) to Sheet1, moves the
# hidden helper list value from Y2 to X2 and repoints the F2 data-validation
# list formula at it, and swaps which sheet/cell is the active
# sheet/selection (Sheet1 becomes the active tab instead of Sheet2).
#
# NOTE: cell values below are intentionally written in a specific order
# (matching how the new values were actually typed in) so that newly
# interned shared strings line up with the canonical workbook.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- New 3rd address column (header + data) --------------------------------
$ws1.Range("J1").Value = "add1"
$ws1.Range("K1").Value = "add2"
$ws1.Range("L1").Value = "add3"
$ws1.Range("J2").Value = "heritage"
$ws1.Range("K2").Value = "city"
$ws1.Range("L2").Value = "sai public school"

# --- Occupation -------------------------------------------------------------
$ws1.Range("M1").Value = "occupation"
$ws1.Range("M2").Value = "Salaried"

# --- Remaining new headers (row 1) ------------------------------------------
$ws1.Range("N1").Value = "Age proof"
$ws1.Range("O1").Value = "HighestEducation"
$ws1.Range("P1").Value = "Annual Income"
$ws1.Range("Q1").Value = "PAN"
$ws1.Range("R1").Value = "Nominee"
$ws1.Range("S1").Value = "relation with nominee"

# --- Remaining new data (row 2) ---------------------------------------------
$ws1.Range("N2").Value = "PAN Card"
$ws1.Range("O2").Value = "Graduate"
$ws1.Range("P2").Value = 1500000
$ws1.Range("Q2").Value = "BUXPG1749Q"
$ws1.Range("R2").Value = "AAC"
$ws1.Range("S2").Value = "Brother"

# --- Move the validation helper value from Y2 to X2 ------------------------
$ws1.Range("X2").Value = $ws1.Range("Y2").Value2
$ws1.Range("Y2").Clear() | Out-Null

# Repoint the "cities" list validation on F2 at its new home ($X$2 instead of
# $Y$2) without disturbing the rest of the validation (prompt title, etc).
$ws1.Range("F2").Validation.Formula1 = "=`$X`$2:`$X`$2"

# --- Sheet2: selection moves from B4 to C5 ----------------------------------
$ws2.Range("C5").Select() | Out-Null

# --- Active sheet/tab switches from Sheet2 to Sheet1, selection S1 --------
$ws1.Activate() | Out-Null
$ws1.Range("S1").Select() | Out-Null
